# Update the "Förändrad" date column (C) for rows 2-5 from 2023-10-05 (45204)
# to 2023-10-08 (45207), matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..5) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45207
    }
}
